$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $value) {
    $rng = $ws.Range($cellRef)
    $origStyle = $rng.Style
    $rng.NumberFormat = "@"
    $rng.Value = $value
    $rng.Style = $origStyle
}

Set-TextValue "D2" "241.39"
Set-TextValue "D3" "23.13"
Set-TextValue "D4" "5.744"
Set-TextValue "D5" "0.05770"
Set-TextValue "D6" "3.420"
Set-TextValue "D7" "6.468"
Set-TextValue "D8" "1.316"
Set-TextValue "D9" "0.8011"
Set-TextValue "D10" "0.1460"
Set-TextValue "D11" "0.07629"
Set-TextValue "D12" "0.03237"
Set-TextValue "D13" "0.02973"
Set-TextValue "D14" "0.09244"
Set-TextValue "D15" "0.001663"
Set-TextValue "D16" "3.258"
Set-TextValue "D17" "0.04758"
Set-TextValue "D18" "0.0005994"
Set-TextValue "D19" "0.006234"
Set-TextValue "D23" "3.694"
Set-TextValue "D25" "0.3322"
Set-TextValue "D26" "0.1276"
Set-TextValue "D27" "0.0006732"
Set-TextValue "D40" "0.04270"
Set-TextValue "D41" "0.007143"
$ws.Range("B42").Value = "CEJI"
$ws.Range("C42").Value = "https://coinranking.com/coin/SbKjCVJCh+ceji-ceji"
Set-TextValue "D42" "0.003446"
$ws.Range("E42").Value = "41CEJICEJI"
$ws.Range("B43").Value = "BKEXToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"
Set-TextValue "D43" "0.1054"
$ws.Range("E43").Value = "42BKEXTokenBKK"
Set-TextValue "D44" "0.009555"
Set-TextValue "D46" "0.00005617"
Set-TextValue "D49" "0.09757"
